$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Gnai2"
$ws.Cells.Item(2,3).Value = "Lhcgr"
$ws.Cells.Item(2,4).Value = "FAPs"
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 151.7260716666667
$ws.Cells.Item(2,8).Value = 455.178215
$ws.Cells.Item(2,9).Value = 0.2700739458961593
$ws.Cells.Item(2,10).Value = 0.2783366498663096
$ws.Cells.Item(2,11).Value = 1
$ws.Cells.Item(2,12).Value = 0.3333333333333333
$ws.Cells.Item(2,13).Value = 0.114644
$ws.Cells.Item(2,14).Value = 0.343932
$ws.Cells.Item(2,15).Value = 0.107929744556041
$ws.Cells.Item(2,16).Value = 0.1536053250846448
$ws.Cells.Item(2,17).Value = 17.39448376015334
$ws.Cells.Item(2,18).Value = 156.55035384138
$ws.Cells.Item(2,19).Value = 0.0291490119918145
$ws.Cells.Item(2,20).Value = 0.04275399158568544

# Row 3
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Gnai2"
$ws.Cells.Item(3,3).Value = "Lhcgr"
$ws.Cells.Item(3,4).Value = "sCs"
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 151.7260716666667
$ws.Cells.Item(3,8).Value = 455.178215
$ws.Cells.Item(3,9).Value = 0.2700739458961593
$ws.Cells.Item(3,10).Value = 0.2783366498663096
$ws.Cells.Item(3,11).Value = 2
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 0.9475655
$ws.Cells.Item(3,14).Value = 1.895131
$ws.Cells.Item(3,15).Value = 0.8920702554439589
$ws.Cells.Item(3,16).Value = 0.8463946749153552
$ws.Cells.Item(3,17).Value = 143.7703909618608
$ws.Cells.Item(3,18).Value = 862.622345771165
$ws.Cells.Item(3,19).Value = 0.2409249339043447
$ws.Cells.Item(3,20).Value = 0.2355826582806241

# Row 4
$ws.Cells.Item(4,1).Value = "FAPs"
$ws.Cells.Item(4,2).Value = "Gnai2"
$ws.Cells.Item(4,3).Value = "Lhcgr"
$ws.Cells.Item(4,4).Value = "FAPs"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 82.248871
$ws.Cells.Item(4,8).Value = 246.746613
$ws.Cells.Item(4,9).Value = 0.146403824289839
$ws.Cells.Item(4,10).Value = 0.150882936320401
$ws.Cells.Item(4,11).Value = 1
$ws.Cells.Item(4,12).Value = 0.3333333333333333
$ws.Cells.Item(4,13).Value = 0.114644
$ws.Cells.Item(4,14).Value = 0.343932
$ws.Cells.Item(4,15).Value = 0.107929744556041
$ws.Cells.Item(4,16).Value = 0.1536053250846448
$ws.Cells.Item(4,17).Value = 9.429339566924
$ws.Cells.Item(4,18).Value = 84.864056102316
$ws.Cells.Item(4,19).Value = 0.01580132735762983
$ws.Cells.Item(4,20).Value = 0.02317642248322095

# Row 5
$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Gnai2"
$ws.Cells.Item(5,3).Value = "Lhcgr"
$ws.Cells.Item(5,4).Value = "sCs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 82.248871
$ws.Cells.Item(5,8).Value = 246.746613
$ws.Cells.Item(5,9).Value = 0.146403824289839
$ws.Cells.Item(5,10).Value = 0.150882936320401
$ws.Cells.Item(5,11).Value = 2
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 0.9475655
$ws.Cells.Item(5,14).Value = 1.895131
$ws.Cells.Item(5,15).Value = 0.8920702554439589
$ws.Cells.Item(5,16).Value = 0.8463946749153552
$ws.Cells.Item(5,17).Value = 77.93619257355049
$ws.Cells.Item(5,18).Value = 467.6171554413029
$ws.Cells.Item(5,19).Value = 0.1306024969322091
$ws.Cells.Item(5,20).Value = 0.12770651383718

# Row 6
$ws.Cells.Item(6,1).Value = "M1"
$ws.Cells.Item(6,2).Value = "Gnai2"
$ws.Cells.Item(6,3).Value = "Lhcgr"
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 123.444321
$ws.Cells.Item(6,8).Value = 370.332963
$ws.Cells.Item(6,9).Value = 0.2197321429647646
$ws.Cells.Item(6,10).Value = 0.2264546783208506
$ws.Cells.Item(6,11).Value = 1
$ws.Cells.Item(6,12).Value = 0.3333333333333333
$ws.Cells.Item(6,13).Value = 0.114644
$ws.Cells.Item(6,14).Value = 0.343932
$ws.Cells.Item(6,15).Value = 0.107929744556041
$ws.Cells.Item(6,16).Value = 0.1536053250846448
$ws.Cells.Item(6,17).Value = 14.152150736724
$ws.Cells.Item(6,18).Value = 127.369356630516
$ws.Cells.Item(6,19).Value = 0.02371563406093852
$ws.Cells.Item(6,20).Value = 0.03478464448041292

# Row 7
$ws.Cells.Item(7,1).Value = "M1"
$ws.Cells.Item(7,2).Value = "Gnai2"
$ws.Cells.Item(7,3).Value = "Lhcgr"
$ws.Cells.Item(7,4).Value = "sCs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 123.444321
$ws.Cells.Item(7,8).Value = 370.332963
$ws.Cells.Item(7,9).Value = 0.2197321429647646
$ws.Cells.Item(7,10).Value = 0.2264546783208506
$ws.Cells.Item(7,11).Value = 2
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 0.9475655
$ws.Cells.Item(7,14).Value = 1.895131
$ws.Cells.Item(7,15).Value = 0.8920702554439589
$ws.Cells.Item(7,16).Value = 0.8463946749153552
$ws.Cells.Item(7,17).Value = 116.9715797505255
$ws.Cells.Item(7,18).Value = 701.8294785031529
$ws.Cells.Item(7,19).Value = 0.196016508903826
$ws.Cells.Item(7,20).Value = 0.1916700338404377

# Row 8
$ws.Cells.Item(8,1).Value = "M2"
$ws.Cells.Item(8,2).Value = "Gnai2"
$ws.Cells.Item(8,3).Value = "Lhcgr"
$ws.Cells.Item(8,4).Value = "FAPs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 154.3429766666667
$ws.Cells.Item(8,8).Value = 463.02893
$ws.Cells.Item(8,9).Value = 0.2747320633285943
$ws.Cells.Item(8,10).Value = 0.2831372788071194
$ws.Cells.Item(8,11).Value = 1
$ws.Cells.Item(8,12).Value = 0.3333333333333333
$ws.Cells.Item(8,13).Value = 0.114644
$ws.Cells.Item(8,14).Value = 0.343932
$ws.Cells.Item(8,15).Value = 0.107929744556041
$ws.Cells.Item(8,16).Value = 0.1536053250846448
$ws.Cells.Item(8,17).Value = 17.69449621697333
$ws.Cells.Item(8,18).Value = 159.25046595276
$ws.Cells.Item(8,19).Value = 0.02965176141640925
$ws.Cells.Item(8,20).Value = 0.04349139375474929

# Row 9
$ws.Cells.Item(9,1).Value = "M2"
$ws.Cells.Item(9,2).Value = "Gnai2"
$ws.Cells.Item(9,3).Value = "Lhcgr"
$ws.Cells.Item(9,4).Value = "sCs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 154.3429766666667
$ws.Cells.Item(9,8).Value = 463.02893
$ws.Cells.Item(9,9).Value = 0.2747320633285943
$ws.Cells.Item(9,10).Value = 0.2831372788071194
$ws.Cells.Item(9,11).Value = 2
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 0.9475655
$ws.Cells.Item(9,14).Value = 1.895131
$ws.Cells.Item(9,15).Value = 0.8920702554439589
$ws.Cells.Item(9,16).Value = 0.8463946749153552
$ws.Cells.Item(9,17).Value = 146.2500798566383
$ws.Cells.Item(9,18).Value = 877.50047913983
$ws.Cells.Item(9,19).Value = 0.245080301912185
$ws.Cells.Item(9,20).Value = 0.2396458850523701

# Row 10
$ws.Cells.Item(10,1).Value = "sCs"
$ws.Cells.Item(10,2).Value = "Gnai2"
$ws.Cells.Item(10,3).Value = "Lhcgr"
$ws.Cells.Item(10,4).Value = "FAPs"
$ws.Cells.Item(10,5).Value = 2
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 50.0323125
$ws.Cells.Item(10,8).Value = 100.064625
$ws.Cells.Item(10,9).Value = 0.08905802352064279
$ws.Cells.Item(10,10).Value = 0.06118845668531954
$ws.Cells.Item(10,11).Value = 1
$ws.Cells.Item(10,12).Value = 0.3333333333333333
$ws.Cells.Item(10,13).Value = 0.114644
$ws.Cells.Item(10,14).Value = 0.343932
$ws.Cells.Item(10,15).Value = 0.107929744556041
$ws.Cells.Item(10,16).Value = 0.1536053250846448
$ws.Cells.Item(10,17).Value = 5.735904434250001
$ws.Cells.Item(10,18).Value = 34.4154266055
$ws.Cells.Item(10,19).Value = 0.009612009729248866
$ws.Cells.Item(10,20).Value = 0.009398872780576216

# Row 11
$ws.Cells.Item(11,1).Value = "sCs"
$ws.Cells.Item(11,2).Value = "Gnai2"
$ws.Cells.Item(11,3).Value = "Lhcgr"
$ws.Cells.Item(11,4).Value = "sCs"
$ws.Cells.Item(11,5).Value = 2
$ws.Cells.Item(11,6).Value = 1
$ws.Cells.Item(11,7).Value = 50.0323125
$ws.Cells.Item(11,8).Value = 100.064625
$ws.Cells.Item(11,9).Value = 0.08905802352064279
$ws.Cells.Item(11,10).Value = 0.06118845668531954
$ws.Cells.Item(11,11).Value = 2
$ws.Cells.Item(11,12).Value = 1
$ws.Cells.Item(11,13).Value = 0.9475655
$ws.Cells.Item(11,14).Value = 1.895131
$ws.Cells.Item(11,15).Value = 0.8920702554439589
$ws.Cells.Item(11,16).Value = 0.8463946749153552
$ws.Cells.Item(11,17).Value = 47.40889321021875
$ws.Cells.Item(11,18).Value = 189.635572840875
$ws.Cells.Item(11,19).Value = 0.07944601379139392
$ws.Cells.Item(11,20).Value = 0.05178958390474333
